# Update loading_percent values for the 380 kV case (rows 2-25, cols B:D,F:I,K:M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","I","K","L","M")

$newValues = @{
    2 = @(11.71488928101106, 8.434182132843857, 6.503814739765549, 32.12301095449332, 42.38729463051282, 17.65249033466434, 28.37936810150292, 9.129534382495686, 11.13554364994146, 14.57061856873929)
    3 = @(11.51764479922584, 8.374148759383495, 6.487727894908961, 32.10506768438444, 42.34391727002451, 17.68638879260773, 28.43534746586974, 8.982034681785848, 11.14338586331154, 14.54863458339762)
    4 = @(11.39757517394168, 8.336297835266707, 6.477562224406262, 32.10141775002383, 42.32806599309615, 17.70995869212606, 28.47460306204254, 8.892315548844078, 11.14987464043311, 14.53758180971919)
    5 = @(11.34897500174641, 8.320628071094303, 6.473347636104686, 32.10178456194593, 42.32432047561178, 17.72025607240428, 28.49182638052129, 8.856016420769707, 11.15294013360899, 14.53369655208956)
    6 = @(11.34092684671063, 8.318011420294187, 6.472643473192534, 32.10195748238745, 42.32386248783968, 17.72200774844783, 28.49476031221553, 8.850006252880508, 11.15347460991257, 14.53308889037766)
    7 = @(11.39691831246402, 8.336087495137342, 6.477505675852401, 32.10141518825557, 42.32800448907133, 17.71009476333156, 28.47483037883242, 8.89182487949962, 11.14991427658813, 14.53752690127438)
    8 = @(11.64670429178034, 8.413690391112583, 6.498328122238599, 32.11529636250108, 42.37010242482308, 17.66360615650043, 28.39765507406214, 9.078530416535552, 11.13790056164057, 14.56253338713263)
    9 = @(12.14182412525159, 8.557851550654071, 6.536852947248571, 32.20085200969435, 42.53800078491177, 17.59433534296529, 28.28514238015255, 9.449218076094359, 11.12759992997982, 14.63078757938668)
    10 = @(12.50493098203343, 8.658670607502762, 6.563724441912759, 32.29903506140134, 42.71296278287974, 17.5568199193856, 28.22624543320505, 9.721499832106833, 11.12808167742614, 14.69237495063467)
    11 = @(12.66918646894317, 8.703386293603057, 6.575632580671923, 32.35129448245518, 42.80362603373723, 17.54266321042823, 28.20462921454118, 9.844774420370319, 11.13003954759935, 14.72281048676537)
    12 = @(12.7311903547364, 8.720150324017581, 6.580095985182243, 32.37216723065859, 42.83953456392188, 17.53772104003015, 28.19718908800529, 9.891324551441324, 11.13102996900242, 14.73467722350083)
    13 = @(12.71784640793019, 8.716547459100731, 6.579136763175551, 32.36762389025319, 42.8317311882662, 17.53876680034713, 28.19875828196678, 9.881305710490309, 11.13080560617766, 14.73210642890067)
    14 = @(12.67429184862792, 8.704768885760599, 6.576000711566481, 32.35299004232029, 42.80654875701205, 17.54224822285497, 28.20400216493406, 9.84860701969596, 11.13011604565641, 14.72377996500585)
    15 = @(12.64758605470279, 8.697532070409782, 6.574073789761963, 32.34416716904592, 42.79132856594847, 17.54443522555367, 28.20731130602996, 9.828559675896921, 11.12972606755217, 14.71872403703009)
    16 = @(12.4941724409237, 8.655724914084933, 6.562939815963379, 32.29577188766602, 42.70725922448665, 17.55780366473583, 28.22776236493051, 9.713427680127619, 11.12798863010589, 14.69043405664184)
    17 = @(12.39977569791052, 8.629780948306065, 6.556028209662482, 32.26802235598117, 42.65851070869237, 17.5667501555573, 28.24163506577525, 9.642613526961329, 11.1273674175098, 14.67369422778083)
    18 = @(12.34539787002355, 8.614750890698703, 6.552023139726915, 32.25277710035285, 42.63151510997974, 17.57216975891336, 28.2501014401035, 9.601830461464804, 11.12717383688424, 14.66429402886777)
    19 = @(12.32697413511225, 8.609643579032063, 6.550662007711703, 32.24773846584423, 42.62255448999813, 17.57405176075586, 28.25305164224709, 9.588014431433658, 11.12713644239161, 14.66115064989079)
    20 = @(12.40983348266289, 8.632553910039153, 6.556767039319354, 32.2709023460481, 42.66359220555528, 17.56576944450006, 28.24010786464616, 9.6501576070317, 11.12741660870909, 14.67545264806034)
    21 = @(12.6870906910523, 8.708233150431296, 6.576923097204268, 32.35725903736276, 42.81390280396118, 17.54121428004176, 28.202441669769, 9.858215343802257, 11.1303118366258, 14.72621643508544)
    22 = @(12.86712394520389, 8.756708117148655, 6.589828176040882, 32.42000700047296, 42.92131856181776, 17.52760646078142, 28.18217015113513, 9.993407612384487, 11.13365495422346, 14.76138056124555)
    23 = @(12.77116406800796, 8.730927562032111, 6.582965195441877, 32.38594326984126, 42.86315467806885, 17.53464583130597, 28.1925915119198, 9.921339767639875, 11.13173827384814, 14.74243320198954)
    24 = @(12.40528668935543, 8.631300610500645, 6.556433112479014, 32.26959809545362, 42.66129164941341, 17.56621196377907, 28.24079678358515, 9.646747146198676, 11.12739385981345, 14.67465696880401)
    25 = @(12.00772169763424, 8.519732205285282, 6.526680712639058, 32.17148112474084, 42.4834774120776, 17.61072757320622, 28.31141264967386, 9.348746501432723, 11.12897009676952, 14.61029242298449)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $row
        $ws.Range($addr).Value = $vals[$i]
    }
}
